$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; this shifts rows 22-82 down to 23-83
# (carrying values/formatting along, matching the diff's net effect).
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new record's data.
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44526
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108002
$ws.Range("J22").Value = "Mango"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 170
$ws.Range("N22").Value = 6000
$ws.Range("O22").Value = 6500
$ws.Range("P22").Value = 6235
$ws.Range("Q22").Value = "$/bandeja 4 kilos"
$ws.Range("R22").Value = "Perú"
$ws.Range("S22").Value = 1559
$ws.Range("T22").Value = 4
